# "commenting and cleaning up" - mark two TODOs as done and add three new
# test/task rows to the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Two items got addressed - flip their "Done?" (column B) flag from N to Y.
$ws.Range("B15").Value = "Y"
$ws.Range("B16").Value = "Y"

# New row: accel-init noise-at-different-ranges idea (Design).
$ws.Range("A25").Value = "in accel init, maybe measure noise at different ranges so we can switch mid flight "
$ws.Range("B25").Value = "N"
$ws.Range("C25").Value = "N"
$ws.Range("D25").Value = "Design"

# New row: altimeter temp/press/alt sanity check (Hardware Test).
$ws.Range("A26").Value = "does the altimeter temp/press/alt match the expected values?"
$ws.Range("B26").Value = "N"
$ws.Range("C26").Value = "N"
$ws.Range("D26").Value = "Hardware Test"

# New row: error/status code logging idea (Design).
$ws.Range("A27").Value = "develop error/status codes that can be written to 'notes' in log file (strings are expensive)"
$ws.Range("B27").Value = "N"
$ws.Range("C27").Value = "N"
$ws.Range("D27").Value = "Design"

# Leave the cursor where the author left it when they saved.
$ws.Range("B17").Select()
